$d = $word.ActiveDocument

# --- Helper: replace a single character inside a Range with another
#     character while forcing Word to keep it as its own run (even
#     though the final formatting is identical to its neighbours).
#     This mirrors what happens when a user selects one character and
#     retypes it: the run gets split around the edited character.
function Replace-CharInRange($start, $len, $newChar) {
    $r = $d.Range($start, $start + $len)
    $r.Text = $newChar
    # Touch formatting so the run is split from its neighbours, then
    # restore it so the resulting rPr is identical to the rest of the
    # paragraph's run (matches target XML exactly).
    $r.Font.Bold = $true
    $r.Font.Bold = $false
}

# --- "Weeks 2-3:" -> "Weeks 2-4:" --------------------------------------
$find1 = $d.Content
$found1 = $find1.Find.Execute("Weeks 2-3:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $digitStart = $find1.Start + "Weeks 2-".Length
    Replace-CharInRange $digitStart 1 "4"
}

# --- "Weeks 4-10:" -> "Weeks 5-10:" ------------------------------------
$find2 = $d.Content
$found2 = $find2.Find.Execute("Weeks 4-10:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $digitStart = $find2.Start + "Weeks ".Length
    Replace-CharInRange $digitStart 1 "5"
}
